$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.987.14"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.501.02"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.13"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.66"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "2.522.32"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "2.942.95"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.44"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "58.903.87"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "2.523.02"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.91"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.76"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.87"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").Value = "0.0₃0773"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.64"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  -7.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.46"
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.39"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.61"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.33"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  -9.09%  "
$ws.Range("E39").Value = "  -4.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.99"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.35"
$ws.Range("E41").Value = "  -7.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.816"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.79"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0928"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.58"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.54"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0513"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0227"
$ws.Range("E51").Value = "  -0.78%  "
